# Updated cryptos list on Tue May 23 19:15:43 UTC 2023 with GitHub Actions
# Applies per-row Price (D) / Volume(1h) (E) updates, and a 3-way row rotation
# of Solana/TRON/WrappedEther (rows 11-13) to match the refreshed rankings feed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.176.98"
$ws.Range("E2").Value = "  +0.68%  "
$ws.Range("D3").Value = "1.849.24"
$ws.Range("E3").Value = "  +1.16%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.47%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.11"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.30%  "
$ws.Range("E6").Value = "  -0.44%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4608"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.26%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07266"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.00%  "
$ws.Range("E10").Value = "  +0.80%  "
$ws.Range("B11").Value = "Solana"
$ws.Range("C11").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.98"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.72%  "
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07825"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.28%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.799.93"
$ws.Range("E13").Value = "  -2.88%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.362"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.36%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.491"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.05%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.10"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.34%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.003"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.39%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008897"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.27%  "
$ws.Range("E19").Value = "  -0.44%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.65"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.13%  "
$ws.Range("D21").Value = "27.208.03"
$ws.Range("E21").Value = "  -0.25%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.051"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.21%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.49"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.62%  "
$ws.Range("D24").Value = "2.121.49"
$ws.Range("E24").Value = "  +3.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.954"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.45%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "152.19"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.55%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.39"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.27%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.049"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.46%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "115.60"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.048"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08804"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.00%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.101"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.61%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7625"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.44%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.167"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.493"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.31%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.716"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +9.83%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.078"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.83%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05228"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.09%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01935"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.97%  "
$ws.Range("E40").Value = "  -0.25%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.039"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5090"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.42%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1624"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.25%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.367"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.21%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4777"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.53%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.34"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.71%  "
$ws.Range("E47").Value = "  -0.49%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "102.68"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.00%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.633"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.08%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06211"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.16%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "65.65"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.00%  "
